$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 4 new rows before the current last data row (row 33) ---
# This pushes the specially-bordered closing row (33) down to row 37,
# and the footer rows (38/39) down to (42/43) automatically.
$ws.Rows("33:36").Insert()

# Copy the formatting of a normal data row (row 32) into the 4 new rows
# so they look like the rest of the table (the closing row keeps its own
# special border formatting since it was shifted down intact).
$ws.Range("B32:J32").Copy($ws.Range("B33:J33"))
$ws.Range("B32:J32").Copy($ws.Range("B34:J34"))
$ws.Range("B32:J32").Copy($ws.Range("B35:J35"))
$ws.Range("B32:J32").Copy($ws.Range("B36:J36"))

# --- Update header summary cells ---
$ws.Range("E11").Value = 989496
$ws.Range("F13").Value = 16

# --- Rewrite the worker/period table (rows 16-37) ---
# Column layout: B=Tipo Doc, C=N Doc, D=Nombre, E=Periodo, F=Valor Mora, G=Salario Basico

# 14 rows for FABIANA ISABEL CASTILLO MENDOZA, periods 1911 .. 2012 (ascending)
$periodos = @("1911","1912","2001","2002","2003","2004","2005","2006","2007","2008","2009","2010","2011","2012")
$r = 16
foreach ($p in $periodos) {
    $ws.Cells.Item($r, 2).Value = "CC"
    $ws.Cells.Item($r, 3).Value = "1238340678"
    $ws.Cells.Item($r, 4).Value = "FABIANA ISABEL CASTILLO MENDOZA"
    $ws.Cells.Item($r, 5).Value = $p
    $ws.Cells.Item($r, 6).Value = 33125
    $ws.Cells.Item($r, 7).Value = 828116
    $r = $r + 1
}

# Remaining four workers, first for period 2507, then period 2508
$workers = @(
    @{ Doc = "45563382"; Nombre = "KELLY DANID OSORIO CASTELLAR"; Mora = 60000; Salario = 1500000 },
    @{ Doc = "1143338851"; Nombre = "UBALDO AMETH BANQUEZ GALVIS"; Mora = 76296; Salario = 1907408 },
    @{ Doc = "1047444358"; Nombre = "DIANA PAOLA GONZALEZ CUADRADO"; Mora = 69637; Salario = 1740915 },
    @{ Doc = "1047490959"; Nombre = "HONEYWELL JOSE SARAVIA SOLANO"; Mora = 56940; Salario = 1423500 }
)

foreach ($periodo in @("2507","2508")) {
    foreach ($w in $workers) {
        $ws.Cells.Item($r, 2).Value = "CC"
        $ws.Cells.Item($r, 3).Value = $w.Doc
        $ws.Cells.Item($r, 4).Value = $w.Nombre
        $ws.Cells.Item($r, 5).Value = $periodo
        $ws.Cells.Item($r, 6).Value = $w.Mora
        $ws.Cells.Item($r, 7).Value = $w.Salario
        $r = $r + 1
    }
}
